$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 56

# New draw-result row appended by the daily auto-update job.
$date       = "2025-11-11"
$game       = "Pick 4"
$phase      = "251111"
$result     = "4-2-0-3"
$insertedAt = "2025-11-11T21:40:24.998+04:00"

# Columns A (date) and C (phase) look like a date / a plain number. A
# direct .Value assignment would make Excel auto-convert them into a real
# date serial / numeric value, but the source sheet stores every column as
# literal text. Format those two cells as Text ("@") before assigning so
# the string is kept verbatim, then switch the cell back to the sheet's
# default "Normal" style so it still renders/serializes like every other
# data row (no per-cell number format lingering on the new row).
$aCell = $ws.Cells.Item($row, 1)
$aCell.NumberFormat = "@"
$aCell.Value = $date
$aCell.Style = "Normal"

$cCell = $ws.Cells.Item($row, 3)
$cCell.NumberFormat = "@"
$cCell.Value = $phase
$cCell.Style = "Normal"

# Columns B, D, E aren't number/date-like, so Excel keeps them as plain
# text on a direct assignment.
$ws.Cells.Item($row, 2).Value = $game
$ws.Cells.Item($row, 4).Value = $result
$ws.Cells.Item($row, 5).Value = $insertedAt
